$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.734.74'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '2.334.13'
$ws.Range('E3').Value = '  -1.13%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.47'
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.660'
$ws.Range('E6').Value = '  -4.42%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '71.58'
$ws.Range('E7').Value = '  -6.22%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.578'
$ws.Range('E9').Value = '  -8.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0979'
$ws.Range('E10').Value = '  -4.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '58.03'
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '32.25'
$ws.Range('E12').Value = '  -3.21%  '
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.08'
$ws.Range('E14').Value = '  -6.25%  '
$ws.Range('D15').Value = '2.682.06'
$ws.Range('E15').Value = '  -0.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.97'
$ws.Range('E16').Value = '  -5.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.888'
$ws.Range('E17').Value = '  -3.55%  '
$ws.Range('D18').Value = '2.338.18'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('D19').Value = '43.626.32'
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('D20').Value = '0.0₂01000'
$ws.Range('E20').Value = '  -3.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '77.57'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.55'
$ws.Range('E22').Value = '  -1.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.79'
$ws.Range('E23').Value = '  -2.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.90'
$ws.Range('E24').Value = '  +7.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  +2.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.49'
$ws.Range('E27').Value = '  -2.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.25'
$ws.Range('E28').Value = '  -8.85%  '
$ws.Range('E29').Value = '  -1.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '175.21'
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.01'
$ws.Range('E31').Value = '  -5.13%  '
$ws.Range('E32').Value = '  -2.56%  '
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0730'
$ws.Range('E34').Value = '  -2.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.02'
$ws.Range('E35').Value = '  -5.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.31'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.72'
$ws.Range('E37').Value = '  -2.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.35'
$ws.Range('E38').Value = '  -1.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.35'
$ws.Range('E39').Value = '  -3.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.51'
$ws.Range('E40').Value = '  +22.18%  '
$ws.Range('E41').Value = '  -3.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '64.60'
$ws.Range('E42').Value = '  +17.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.14'
$ws.Range('E43').Value = '  +2.72%  '
$ws.Range('E44').Value = '  +2.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.62'
$ws.Range('E45').Value = '  -1.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.194'
$ws.Range('E46').Value = '  -3.54%  '
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.21'
$ws.Range('E48').Value = '  -3.99%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.94'
$ws.Range('E49').Value = '  +4.35%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.41'
$ws.Range('E50').Value = '  -4.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '97.29'
$ws.Range('E51').Value = '  -4.59%  '
